$d = $word.ActiveDocument

function Set-ParaText($para, [string]$newText) {
    # Replace a paragraph's content (all runs) with a single run containing $newText,
    # dropping any proofErr spell-check markers and merging multiple runs into one.
    $full = $para.Range
    $body = $d.Range($full.Start, $full.End - 1)   # exclude the paragraph mark
    if ($body.Start -lt $body.End) {
        $body.Delete()
    }
    $collapsed = $d.Range($full.Start, $full.Start)
    $collapsed.InsertBefore($newText)
}

# ---------------------------------------------------------------------------
# Change 1: "Van <spellStart>Quy</spellEnd> Tran: " -> single run "Van Quy Tran: "
# ---------------------------------------------------------------------------
Set-ParaText $d.Paragraphs.Item(1) "Van Quy Tran: "

# ---------------------------------------------------------------------------
# Change 2: "Thanh <spellStart>Huy</spellEnd> Le" -> single run "Thanh Huy Le"
# ---------------------------------------------------------------------------
Set-ParaText $d.Paragraphs.Item(18) "Thanh Huy Le"

# ---------------------------------------------------------------------------
# Change 3: "Rest Controller" paragraph gains a second run " for Comment"
# ---------------------------------------------------------------------------
$pRestController = $d.Paragraphs.Item(19)
$pRestController.Range.InsertAfter(" for Comment")

# ---------------------------------------------------------------------------
# Change 4: restructure the tail of the document.
#   - 5 new "List Paragraph" bullet items are added right after "Footer"
#   - the old content (blank line, Ferdis Fernando, Blog post/list/detail,
#     trailing blank) shifts down, unchanged in substance
#   - "Ferdis Fernando" loses its proofErr markup / run split
#   - the _GoBack bookmark moves from the trailing blank paragraph to the
#     end of the "Footer" paragraph
# ---------------------------------------------------------------------------
$pFooter = $d.Paragraphs.Item(21)

$newItems = @(
    "Individually Exception Handler",
    "Comment Validation",
    "Tiles for Comments section",
    "Java Configuration for Comment",
    "Internationalization"
)

$insertAfter = $pFooter
foreach ($itemText in $newItems) {
    $insertAfter.Range.InsertParagraphAfter()
    $insertAfter = $d.Paragraphs.Item($insertAfter.Index + 1)
    $insertAfter.Range.InsertBefore($itemText)
}

# Ferdis Fernando now sits 5 paragraphs further down than before (index 23 -> 28)
Set-ParaText $d.Paragraphs.Item(28) "Ferdis Fernando"

# ---- Move the _GoBack bookmark from the final blank paragraph up to the
#      end of the "Footer" paragraph. --------------------------------------
# Collapse every paragraph between "Footer" and the bookmark's paragraph
# (inclusive) into a single run of text by deleting the paragraph marks
# that separate them; the zero-width bookmark travels with the text that
# follows it, ending up glued to the very end of the combined blob.
$lastParaIndex = $d.Paragraphs.Count
$mergeCount = $lastParaIndex - $pFooter.Index
for ($k = 0; $k -lt $mergeCount; $k++) {
    $mark = $d.Range($pFooter.Range.End - 1, $pFooter.Range.End)
    $mark.Delete()
}

$blob = $d.Paragraphs.Item($pFooter.Index)
$blobStart = $blob.Range.Start

# Re-derive the pieces of the combined blob in document order so the split
# offsets always match the live text, regardless of exact wording above.
$pieces = @("Footer") + $newItems + @("", "Ferdis Fernando", "Blog post", "Blog list by category", "Blog detail")

$cumulative = 0
$cuts = @()
for ($i = 0; $i -lt ($pieces.Length - 1); $i++) {
    $cumulative += $pieces[$i].Length
    $cuts += $cumulative
}

# Apply the cuts from the end of the blob backwards so earlier offsets stay valid.
for ($i = $cuts.Length - 1; $i -ge 0; $i--) {
    $pt = $blobStart + $cuts[$i]
    $insPoint = $d.Range($pt, $pt)
    $insPoint.InsertParagraphAfter()
}

# The final (empty) ListParagraph that used to carry the bookmark keeps the
# "List Paragraph" style but no longer hosts the bookmark.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
Write-Host "Final paragraph count:" $d.Paragraphs.Count
Write-Host "Footer paragraph text:" $d.Paragraphs.Item($pFooter.Index).Range.Text
Write-Host "Last paragraph text: [" $lastPara.Range.Text "]"
